$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Myoc"
$ws.Range("C2").Value = "Fzd1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2144083333333333
$ws.Range("H2").Value = 0.643225
$ws.Range("I2").Value = 0.008611346839948651
$ws.Range("J2").Value = 0.008611346839948651
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.6559766666666667
$ws.Range("N2").Value = 1.96793
$ws.Range("O2").Value = 0.030799191223283
$ws.Range("P2").Value = 0.030799191223283
$ws.Range("Q2").Value = 0.1406468638055555
$ws.Range("R2").Value = 1.26582177425
$ws.Range("S2").Value = 0.0002652225180135923
$ws.Range("T2").Value = 0.0002652225180135923
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Myoc"
$ws.Range("C3").Value = "Fzd1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2144083333333333
$ws.Range("H3").Value = 0.643225
$ws.Range("I3").Value = 0.008611346839948651
$ws.Range("J3").Value = 0.008611346839948651
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 16.81477433333333
$ws.Range("N3").Value = 50.444323
$ws.Range("O3").Value = 0.7894815111340611
$ws.Range("P3").Value = 0.789481511134061
$ws.Range("Q3").Value = 3.605227740186111
$ws.Range("R3").Value = 32.447049661675
$ws.Range("S3").Value = 0.006798499116102183
$ws.Range("T3").Value = 0.006798499116102182
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Myoc"
$ws.Range("C4").Value = "Fzd1"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.2144083333333333
$ws.Range("H4").Value = 0.643225
$ws.Range("I4").Value = 0.008611346839948651
$ws.Range("J4").Value = 0.008611346839948651
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.5954103333333334
$ws.Range("N4").Value = 1.786231
$ws.Range("O4").Value = 0.02795550153610953
$ws.Range("P4").Value = 0.02795550153610953
$ws.Range("Q4").Value = 0.1276609372194445
$ws.Range("R4").Value = 1.148948434975
$ws.Range("S4").Value = 0.0002407345198121564
$ws.Range("T4").Value = 0.0002407345198121564
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Myoc"
$ws.Range("C5").Value = "Fzd1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.2144083333333333
$ws.Range("H5").Value = 0.643225
$ws.Range("I5").Value = 0.008611346839948651
$ws.Range("J5").Value = 0.008611346839948651
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.232341666666667
$ws.Range("N5").Value = 9.697025
$ws.Range("O5").Value = 0.1517637961065464
$ws.Range("P5").Value = 0.1517637961065464
$ws.Range("Q5").Value = 0.6930409895138889
$ws.Range("R5").Value = 6.237368905625001
$ws.Range("S5").Value = 0.00130689068602072
$ws.Range("T5").Value = 0.00130689068602072
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Myoc"
$ws.Range("C6").Value = "Fzd1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 23.971258
$ws.Range("H6").Value = 71.91377399999999
$ws.Range("I6").Value = 0.9627648963950115
$ws.Range("J6").Value = 0.9627648963950115
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.6559766666666667
$ws.Range("N6").Value = 1.96793
$ws.Range("O6").Value = 0.030799191223283
$ws.Range("P6").Value = 0.030799191223283
$ws.Range("Q6").Value = 15.72458591864666
$ws.Range("R6").Value = 141.52127326782
$ws.Range("S6").Value = 0.02965238014713421
$ws.Range("T6").Value = 0.02965238014713421
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Myoc"
$ws.Range("C7").Value = "Fzd1"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 23.971258
$ws.Range("H7").Value = 71.91377399999999
$ws.Range("I7").Value = 0.9627648963950115
$ws.Range("J7").Value = 0.9627648963950115
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 16.81477433333333
$ws.Range("N7").Value = 50.444323
$ws.Range("O7").Value = 0.7894815111340611
$ws.Range("P7").Value = 0.789481511134061
$ws.Range("Q7").Value = 403.0712937561112
$ws.Range("R7").Value = 3627.641643805001
$ws.Range("S7").Value = 0.7600850852727614
$ws.Range("T7").Value = 0.7600850852727613
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Myoc"
$ws.Range("C8").Value = "Fzd1"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 23.971258
$ws.Range("H8").Value = 71.91377399999999
$ws.Range("I8").Value = 0.9627648963950115
$ws.Range("J8").Value = 0.9627648963950115
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.5954103333333334
$ws.Range("N8").Value = 1.786231
$ws.Range("O8").Value = 0.02795550153610953
$ws.Range("P8").Value = 0.02795550153610953
$ws.Range("Q8").Value = 14.27273471619933
$ws.Range("R8").Value = 128.454612445794
$ws.Range("S8").Value = 0.02691457554008307
$ws.Range("T8").Value = 0.02691457554008307
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Myoc"
$ws.Range("C9").Value = "Fzd1"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 23.971258
$ws.Range("H9").Value = 71.91377399999999
$ws.Range("I9").Value = 0.9627648963950115
$ws.Range("J9").Value = 0.9627648963950115
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.232341666666667
$ws.Range("N9").Value = 9.697025
$ws.Range("O9").Value = 0.1517637961065464
$ws.Range("P9").Value = 0.1517637961065464
$ws.Range("Q9").Value = 77.48329603581665
$ws.Range("R9").Value = 697.3496643223499
$ws.Range("S9").Value = 0.1461128554350328
$ws.Range("T9").Value = 0.1461128554350328
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Myoc"
$ws.Range("C10").Value = "Fzd1"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.7126843333333334
$ws.Range("H10").Value = 2.138053
$ws.Range("I10").Value = 0.02862375676503981
$ws.Range("J10").Value = 0.02862375676503981
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.6559766666666667
$ws.Range("N10").Value = 1.96793
$ws.Range("O10").Value = 0.030799191223283
$ws.Range("P10").Value = 0.030799191223283
$ws.Range("Q10").Value = 0.4675042933655555
$ws.Range("R10").Value = 4.20753864029
$ws.Range("S10").Value = 0.0008815885581352016
$ws.Range("T10").Value = 0.0008815885581352016
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Myoc"
$ws.Range("C11").Value = "Fzd1"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.7126843333333334
$ws.Range("H11").Value = 2.138053
$ws.Range("I11").Value = 0.02862375676503981
$ws.Range("J11").Value = 0.02862375676503981
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 16.81477433333333
$ws.Range("N11").Value = 50.444323
$ws.Range("O11").Value = 0.7894815111340611
$ws.Range("P11").Value = 0.789481511134061
$ws.Range("Q11").Value = 11.98362623590211
$ws.Range("R11").Value = 107.852636123119
$ws.Range("S11").Value = 0.02259792674519744
$ws.Range("T11").Value = 0.02259792674519743
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Myoc"
$ws.Range("C12").Value = "Fzd1"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.7126843333333334
$ws.Range("H12").Value = 2.138053
$ws.Range("I12").Value = 0.02862375676503981
$ws.Range("J12").Value = 0.02862375676503981
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.5954103333333334
$ws.Range("N12").Value = 1.786231
$ws.Range("O12").Value = 0.02795550153610953
$ws.Range("P12").Value = 0.02795550153610953
$ws.Range("Q12").Value = 0.4243396164714445
$ws.Range("R12").Value = 3.819056548243001
$ws.Range("S12").Value = 0.0008001914762142959
$ws.Range("T12").Value = 0.0008001914762142959
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Myoc"
$ws.Range("C13").Value = "Fzd1"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.7126843333333334
$ws.Range("H13").Value = 2.138053
$ws.Range("I13").Value = 0.02862375676503981
$ws.Range("J13").Value = 0.02862375676503981
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 3.232341666666667
$ws.Range("N13").Value = 9.697025
$ws.Range("O13").Value = 0.1517637961065464
$ws.Range("P13").Value = 0.1517637961065464
$ws.Range("Q13").Value = 2.303639265813889
$ws.Range("R13").Value = 20.732753392325
$ws.Range("S13").Value = 0.00434404998549288
$ws.Range("T13").Value = 0.00434404998549288
